$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their values as text (avoid Excel auto-converting
# numeric-looking strings such as "1.013" or "161.00" into numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.372.39'
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").Value = '1.874.84'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +1.07%  '

$ws.Range("D5").Value = '314.98'
$ws.Range("E5").Value = '  +0.61%  '

$ws.Range("D6").Value = '1.010'
$ws.Range("E6").Value = '  +0.83%  '

$ws.Range("D7").Value = '0.5135'
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").Value = '0.3935'
$ws.Range("E8").Value = '  +1.17%  '

$ws.Range("D9").Value = '0.08334'
$ws.Range("E9").Value = '  -0.63%  '

$ws.Range("D10").Value = '1.118'
$ws.Range("E10").Value = '  +0.41%  '

$ws.Range("D11").Value = '41.82'
$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("D12").Value = '6.257'
$ws.Range("E12").Value = '  +1.06%  '

$ws.Range("D13").Value = '20.34'
$ws.Range("E13").Value = '  -1.12%  '

$ws.Range("D14").Value = '1.848.72'
$ws.Range("E14").Value = '  -1.58%  '

$ws.Range("D15").Value = '7.246'
$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D16").Value = '1.015'
$ws.Range("E16").Value = '  +1.24%  '

$ws.Range("D17").Value = '0.00001105'
$ws.Range("E17").Value = '  -0.42%  '

$ws.Range("D18").Value = '91.44'
$ws.Range("E18").Value = '  +0.62%  '

$ws.Range("D19").Value = '0.06719'
$ws.Range("E19").Value = '  +1.07%  '

$ws.Range("D20").Value = '17.73'
$ws.Range("E20").Value = '  +0.34%  '

$ws.Range("D21").Value = '1.010'
$ws.Range("E21").Value = '  +0.75%  '

$ws.Range("D22").Value = '5.970'
$ws.Range("E22").Value = '  -0.92%  '

$ws.Range("D23").Value = '28.359.19'
$ws.Range("E23").Value = '  +0.88%  '

$ws.Range("D24").Value = '11.13'
$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("D25").Value = '2.250'
$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("D26").Value = '2.058.68'
$ws.Range("E26").Value = '  -1.19%  '

$ws.Range("D27").Value = '161.00'
$ws.Range("E27").Value = '  +1.91%  '

$ws.Range("D28").Value = '20.73'
$ws.Range("E28").Value = '  +0.73%  '

$ws.Range("D29").Value = '2.415'
$ws.Range("E29").Value = '  -2.49%  '

$ws.Range("D30").Value = '126.94'
$ws.Range("E30").Value = '  +1.66%  '

$ws.Range("D31").Value = '0.1059'
$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("D32").Value = '1.039'
$ws.Range("E32").Value = '  +0.19%  '

$ws.Range("D33").Value = '5.875'
$ws.Range("E33").Value = '  -0.21%  '

$ws.Range("D34").Value = '3.630'
$ws.Range("E34").Value = '  +0.93%  '

$ws.Range("D35").Value = '0.02448'
$ws.Range("E35").Value = '  +0.39%  '

$ws.Range("D36").Value = '0.06506'
$ws.Range("E36").Value = '  -0.29%  '

$ws.Range("D37").Value = '9.145'
$ws.Range("E37").Value = '  -5.06%  '

$ws.Range("D38").Value = '0.2180'
$ws.Range("E38").Value = '  -0.26%  '

$ws.Range("D39").Value = '1.252'
$ws.Range("E39").Value = '  +2.08%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '1.190'
$ws.Range("E40").Value = '  -1.47%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6459'
$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("D42").Value = '4.976'
$ws.Range("E42").Value = '  -0.55%  '

$ws.Range("D43").Value = '11.11'
$ws.Range("E43").Value = '  -1.83%  '

$ws.Range("D44").Value = '0.6039'
$ws.Range("E44").Value = '  -0.78%  '

$ws.Range("D45").Value = '12.97'
$ws.Range("E45").Value = '  -0.45%  '

$ws.Range("D46").Value = '3.688'
$ws.Range("E46").Value = '  +0.44%  '

$ws.Range("D47").Value = '1.280'
$ws.Range("E47").Value = '  +0.19%  '

$ws.Range("D48").Value = '2.011'
$ws.Range("E48").Value = '  +0.28%  '

$ws.Range("D49").Value = '1.210'
$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("D50").Value = '121.34'
$ws.Range("E50").Value = '  +0.00%  '

$ws.Range("D51").Value = '0.06887'
$ws.Range("E51").Value = '  +0.25%  '
